$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextCell $ws "D2" "35.373.99"
Set-TextCell $ws "E2" "  -0.06%  "
Set-TextCell $ws "D3" "1.915.42"
Set-TextCell $ws "E3" "  +0.73%  "
Set-TextCell $ws "E4" "  -0.84%  "
Set-TextCell $ws "D5" "0.730"
Set-TextCell $ws "E5" "  +12.89%  "
Set-TextCell $ws "D6" "251.69"
Set-TextCell $ws "E6" "  +3.36%  "
Set-TextCell $ws "E7" "  -0.59%  "
Set-TextCell $ws "D8" "40.60"
Set-TextCell $ws "E8" "  -1.84%  "
Set-TextCell $ws "D9" "0.357"
Set-TextCell $ws "E9" "  +4.26%  "
Set-TextCell $ws "D10" "52.75"
Set-TextCell $ws "E10" "  +5.75%  "
Set-TextCell $ws "D11" "0.0735"
Set-TextCell $ws "E11" "  +3.87%  "
Set-TextCell $ws "E12" "  +0.00%  "
Set-TextCell $ws "D13" "2.190.59"
Set-TextCell $ws "E13" "  +0.39%  "
Set-TextCell $ws "D14" "12.58"
Set-TextCell $ws "E14" "  +4.07%  "
Set-TextCell $ws "E15" "  +3.20%  "
Set-TextCell $ws "D16" "1.912.70"
Set-TextCell $ws "E16" "  -0.19%  "
Set-TextCell $ws "D17" "4.90"
Set-TextCell $ws "E17" "  +1.56%  "
Set-TextCell $ws "D18" "35.355.62"
Set-TextCell $ws "E18" "  -0.20%  "
Set-TextCell $ws "D19" "73.15"
Set-TextCell $ws "E19" "  +1.86%  "
Set-TextCell $ws "D20" "0.0₃0831"
Set-TextCell $ws "E20" "  +2.38%  "
Set-TextCell $ws "D21" "13.07"
Set-TextCell $ws "E21" "  +4.77%  "
Set-TextCell $ws "D22" "241.98"
Set-TextCell $ws "E22" "  +0.34%  "
Set-TextCell $ws "D23" "5.07"
Set-TextCell $ws "E23" "  +6.96%  "
Set-TextCell $ws "E24" "  -0.65%  "
Set-TextCell $ws "D25" "2.34"
Set-TextCell $ws "E25" "  +2.30%  "
Set-TextCell $ws "D26" "2.32"
Set-TextCell $ws "E26" "  +5.27%  "
Set-TextCell $ws "D27" "167.84"
Set-TextCell $ws "E27" "  -1.50%  "
Set-TextCell $ws "D28" "8.75"
Set-TextCell $ws "E28" "  +5.45%  "
Set-TextCell $ws "D29" "0.135"
Set-TextCell $ws "E29" "  +7.55%  "
Set-TextCell $ws "D30" "18.76"
Set-TextCell $ws "E30" "  +3.86%  "
Set-TextCell $ws "D31" "4.127.65"
Set-TextCell $ws "E31" "  +19.43%  "
Set-TextCell $ws "D32" "4.37"
Set-TextCell $ws "E32" "  +6.28%  "
Set-TextCell $ws "D33" "1.99"
Set-TextCell $ws "E33" "  +13.93%  "
Set-TextCell $ws "D34" "1.64"
Set-TextCell $ws "E34" "  +24.37%  "
Set-TextCell $ws "D35" "0.0579"
Set-TextCell $ws "E35" "  +3.35%  "
Set-TextCell $ws "D36" "4.25"
Set-TextCell $ws "E36" "  +3.75%  "
Set-TextCell $ws "E37" "  -1.00%  "
Set-TextCell $ws "D38" "0.909"
Set-TextCell $ws "E38" "  -0.43%  "
Set-TextCell $ws "D39" "2.04"
Set-TextCell $ws "E39" "  +0.27%  "
Set-TextCell $ws "D40" "17.46"
Set-TextCell $ws "E40" "  +12.13%  "
Set-TextCell $ws "D41" "98.90"
Set-TextCell $ws "E41" "  +10.51%  "
Set-TextCell $ws "D42" "1.14"
Set-TextCell $ws "E42" "  +4.91%  "
Set-TextCell $ws "D43" "0.0210"
Set-TextCell $ws "E43" "  +1.02%  "
Set-TextCell $ws "D44" "0.0649"
Set-TextCell $ws "E44" "  +2.39%  "
Set-TextCell $ws "B45" "Maker"
Set-TextCell $ws "C45" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws "D45" "1.349.42"
Set-TextCell $ws "E45" "  +1.09%  "
Set-TextCell $ws "B46" "RenderToken"
Set-TextCell $ws "C46" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws "D46" "2.46"
Set-TextCell $ws "E46" "  +5.10%  "
Set-TextCell $ws "E47" "  +0.27%  "
Set-TextCell $ws "D48" "6.73"
Set-TextCell $ws "E48" "  +3.46%  "
Set-TextCell $ws "E49" "  -0.10%  "
Set-TextCell $ws "D50" "45.27"
Set-TextCell $ws "E50" "  -2.84%  "
Set-TextCell $ws "B51" "Gas"
Set-TextCell $ws "C51" "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
Set-TextCell $ws "D51" "11.73"
Set-TextCell $ws "E51" "  -3.85%  "
